# 1. Remove the stray empty B6 cell on the "ODI Batting" sheet.
$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B6").Value = $null

# 2. Add a new worksheet "ODI Batting Extra" as the last (3rd) sheet.
$playerInfo = $wb.Worksheets.Item("Player Info")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Copy the bold/centered header style used on the other sheets' header rows.
$playerInfo.Range("A1:D1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# 3. Populate header row.
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

function Set-TextCell($sheet, $rowNum, $colNum, $text) {
    $cell = $sheet.Cells.Item($rowNum, $colNum)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# 4. Populate data rows (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH).
# Row 2 - match 4485
Set-TextCell $extra 2 1 "4485"
$extra.Range("B2").Value = 3
Set-TextCell $extra 2 3 "5"
Set-TextCell $extra 2 4 "1"
Set-TextCell $extra 2 5 "20.44%"
$extra.Range("F2").Value = "NO"

# Row 3 - match 4621
Set-TextCell $extra 3 1 "4621"
$extra.Range("B3").Value = 5
Set-TextCell $extra 3 3 "0"
Set-TextCell $extra 3 4 "1"
Set-TextCell $extra 3 5 "3.90%"
$extra.Range("F3").Value = "NO"

# Row 4 - match 4623 (only MATCH_CODE and MAN_OF_MATCH known)
Set-TextCell $extra 4 1 "4623"
$extra.Range("F4").Value = "NO"

# Row 5 - match 4624
Set-TextCell $extra 5 1 "4624"
$extra.Range("B5").Value = 5
Set-TextCell $extra 5 3 "0"
Set-TextCell $extra 5 4 "0"
Set-TextCell $extra 5 5 "2.67%"
$extra.Range("F5").Value = "NO"

# Row 6 - match 4637 (only MATCH_CODE, BATTING_POSITION and MAN_OF_MATCH known)
Set-TextCell $extra 6 1 "4637"
$extra.Range("B6").Value = 6
$extra.Range("F6").Value = "NO"

# Row 7 - match 4640
Set-TextCell $extra 7 1 "4640"
$extra.Range("B7").Value = 6
Set-TextCell $extra 7 3 "3"
Set-TextCell $extra 7 4 "4"
Set-TextCell $extra 7 5 "25.75%"
$extra.Range("F7").Value = "YES"

# Row 8 - match 4643 (only MATCH_CODE and MAN_OF_MATCH known)
Set-TextCell $extra 8 1 "4643"
$extra.Range("F8").Value = "NO"

# Row 9 - match 4656 (only MATCH_CODE and MAN_OF_MATCH known)
Set-TextCell $extra 9 1 "4656"
$extra.Range("F9").Value = "NO"

# Row 10 - match 4657
Set-TextCell $extra 10 1 "4657"
$extra.Range("B10").Value = 5
Set-TextCell $extra 10 3 "1"
Set-TextCell $extra 10 4 "1"
Set-TextCell $extra 10 5 "10.64%"
$extra.Range("F10").Value = "NO"

# Row 11 - match 4658
Set-TextCell $extra 11 1 "4658"
$extra.Range("B11").Value = 5
Set-TextCell $extra 11 3 "0"
Set-TextCell $extra 11 4 "0"
Set-TextCell $extra 11 5 "1.90%"
$extra.Range("F11").Value = "NO"

# Row 12 - match 4669
Set-TextCell $extra 12 1 "4669"
$extra.Range("B12").Value = 6
Set-TextCell $extra 12 3 "4"
Set-TextCell $extra 12 4 "0"
Set-TextCell $extra 12 5 "11.76%"
$extra.Range("F12").Value = "NO"

Write-Host "Done applying edits"
